{"js": "// 1) Curso (semestre ideal): swap EA/EB semester numbers.\nconst found = context.document.body.search(\"EA (5), EB (6)\", { matchCase: true });\nfound.load(\"text\");\nawait context.sync();\n\nif (found.items.length > 0) {\n  found.items[0].insertText(\"EA (6), EB (5)\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Reorder the Requisitos bullet list: move the \"LOQ4083 - ...\" weak\n//    requirement (currently last) so it becomes the first entry, ahead of\n//    \"LOB1006 - ...\" and \"LOB1019 - ...\".\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"LOB1006\") !== -1 && t.indexOf(\"LOQ4083\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (target) {\n  const loqText = \"LOQ4083 -  Fen\u00f4menos de Transporte I  (Requisito fraco)\";\n\n  // Insert a fresh copy of the LOQ4083 run (with its own line break) at the\n  // very start of the paragraph, as its own run.\n  const startRange = target.getRange(\"Start\");\n  startRange.insertText(loqText + \"\\u000b\", \"Before\");\n  await context.sync();\n\n  // Remove the original (now duplicated) LOQ4083 run, which sits at the\n  // tail end of the paragraph.\n  const matches = target.search(\"LOQ4083\", { matchCase: true });\n  matches.load(\"text\");\n  await context.sync();\n\n  const lastMatch = matches.items[matches.items.length - 1];\n  const trailingRange = lastMatch.expandTo(target.getRange(\"End\"));\n  trailingRange.delete();\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Curso (semestre ideal): swap EA/EB semester numbers\n$range = $d.Content\n$range.Find.Execute(\"EA (5), EB (6)\", $false, $false, $false, $false, $false, $true, 1, $false, \"EA (6), EB (5)\", 2)\n\n# 2) Reorder the Requisitos bullet list: move the \"LOQ4083 - ...\" weak\n#    requirement (currently last) so it becomes the first entry, ahead of\n#    \"LOB1006 - ...\" and \"LOB1019 - ...\".\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    if ($para.Range.Text -like \"LOB1006*\" -and $para.Range.Text -like \"*LOQ4083*\") {\n        $target = $para\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $fullText = $target.Range.Text\n    $loqText = \"LOQ4083 -  Fen\u00f4menos de Transporte I  (Requisito fraco)\"\n\n    # Insert a fresh copy of the LOQ4083 run (with its own line break) at\n    # the very start of the paragraph, as its own run.\n    $startPoint = $d.Range($target.Range.Start, $target.Range.Start)\n    $startPoint.InsertBefore($loqText + \"`v\")\n\n    # Remove the original (now duplicated) LOQ4083 run, which sits at the\n    # tail end of the paragraph.\n    $lastIdx = $target.Range.Text.LastIndexOf(\"LOQ4083\")\n    $delStart = $target.Range.Start + $lastIdx\n    $delEnd = $target.Range.End\n    $d.Range($delStart, $delEnd).Delete()\n}\n"}
